$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.904.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.324.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.537"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.336.83"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.155"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.60%  "
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.750.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "56.831.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.92%  "
$ws.Range("E17").Value = "  +4.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.335.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.22%  "
$ws.Range("E20").Value = "  +4.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.93%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.994"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.158"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.49%  "
$ws.Range("E27").Value = "  +5.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "171.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  +12.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0735"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.18%  "
$ws.Range("E32").Value = "  +4.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.35%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("B36").Value = "SuiNetwork"
$ws.Range("C36").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.947"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.94%  "
$ws.Range("E38").Value = "  +9.05%  "
$ws.Range("E39").Value = "  +9.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +13.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.17%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "276.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +14.97%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.43%  "
$ws.Range("E46").Value = "  +3.71%  "
$ws.Range("E47").Value = "  +4.08%  "
$ws.Range("E48").Value = "  +3.28%  "
$ws.Range("B49").Value = "Polygon"
$ws.Range("C49").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.382"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0216"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.13%  "
